# project-evaluation.xlsx update: lecture,etude registry and schedule
$wb = $excel.ActiveWorkbook

$wsMember1 = $wb.Worksheets.Item("Member 1")
$wsMember2 = $wb.Worksheets.Item("Member 2")

# --- Member 1 sheet: "time" data type for the Lectures table changes from
#     TIME(0) to TIME ---
$wsMember1.Range("E8").Value = "TIME"

# --- Member 2 sheet: Lectures (D/E) and Etudes (G/H) tables gain a new
#     "profile"/BYTEA non-key column at the top, pushing the existing rows
#     down by one (joindate/DATE, experienceyear/INT, subject/VARCHAR (255)).
#     The old "email" entry in the Etudes column is replaced along the way
#     by the shifted "experienceyear" value. ---
$wsMember2.Range("D8").Value = "profile"
$wsMember2.Range("E8").Value = "BYTEA"
$wsMember2.Range("G8").Value = "profile"
$wsMember2.Range("H8").Value = "BYTEA"

$wsMember2.Range("D9").Value = "joindate"
$wsMember2.Range("E9").Value = "DATE"
$wsMember2.Range("G9").Value = "joindate"
$wsMember2.Range("H9").Value = "DATE"

$wsMember2.Range("D10").Value = "experienceyear"
$wsMember2.Range("E10").Value = "INT"
$wsMember2.Range("G10").Value = "experienceyear"
$wsMember2.Range("H10").Value = "INT"

$wsMember2.Range("D11").Value = "subject"
$wsMember2.Range("E11").Value = "VARCHAR (255)"

# --- Member 2 sheet: UPDATE functionality status moves from Done to
#     Partial for all three main tables ---
$wsMember2.Range("B20").Value = "Partial"
$wsMember2.Range("E20").Value = "Partial"
$wsMember2.Range("H20").Value = "Partial"

# --- Member 2 sheet: BLOBS / FILE UPLOADS status updated ---
$wsMember2.Range("B47").Value = "Yes/No"
$wsMember2.Range("B48").Value = "Yes/No"

# --- Selection / active tab bookkeeping: Member 1 becomes the active tab,
#     with cell E10 selected; Member 2 keeps cell G41 selected. ---
$wsMember2.Range("G41").Select()
$wsMember1.Activate()
$wsMember1.Range("E10").Select()
